$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 1000.2  # H38
$ws.Cells.Item(38, 10).Value = 4000  # J38
$ws.Cells.Item(38, 12).Value = 12000  # L38
$ws.Cells.Item(38, 14).Value = -12744  # N38
$ws.Cells.Item(86, 8).Value = 4016.7083  # H86
$ws.Cells.Item(86, 9).Value = 3508.1667  # I86
$ws.Cells.Item(86, 11).Value = 3508.1667  # K86
$ws.Cells.Item(86, 13).Value = -2385.1667  # M86
$ws.Cells.Item(89, 8).Value = 4016.7083  # H89
$ws.Cells.Item(89, 9).Value = 3508.1667  # I89
$ws.Cells.Item(89, 11).Value = 17540.8335  # K89
$ws.Cells.Item(89, 13).Value = -11924.8335  # M89
$ws.Cells.Item(98, 8).Value = 3682.8147  # H98
$ws.Cells.Item(98, 9).Value = 4056.2273  # I98
$ws.Cells.Item(98, 10).Value = 2039.8  # J98
$ws.Cells.Item(98, 11).Value = 4056.2273  # K98
$ws.Cells.Item(98, 12).Value = 2039.8  # L98
$ws.Cells.Item(98, 13).Value = -2558.2273  # M98
$ws.Cells.Item(98, 14).Value = -5035.8  # N98
$ws.Cells.Item(122, 8).Value = 3682.8147  # H122
$ws.Cells.Item(122, 9).Value = 4056.2273  # I122
$ws.Cells.Item(122, 10).Value = 2039.8  # J122
$ws.Cells.Item(122, 11).Value = 12168.6819  # K122
$ws.Cells.Item(122, 12).Value = 6119.4  # L122
$ws.Cells.Item(122, 13).Value = -9718.6819  # M122
$ws.Cells.Item(122, 14).Value = -11019.4  # N122
$ws.Cells.Item(123, 8).Value = 28999  # H123
$ws.Cells.Item(123, 10).Value = 28999  # J123
$ws.Cells.Item(123, 12).Value = 28999  # L123
$ws.Cells.Item(123, 14).Value = -38799  # N123
$ws.Cells.Item(125, 8).Value = 1729.1  # H125
$ws.Cells.Item(125, 9).Value = 1800.3334  # I125
$ws.Cells.Item(125, 10).Value = 1698.5714  # J125
$ws.Cells.Item(125, 11).Value = 16203.0006  # K125
$ws.Cells.Item(125, 12).Value = 15287.1426  # L125
$ws.Cells.Item(125, 13).Value = -13743.0006  # M125
$ws.Cells.Item(125, 14).Value = -20207.1426  # N125
$ws.Cells.Item(137, 8).Value = 1023.5543  # H137
$ws.Cells.Item(137, 9).Value = 859.7843  # I137
$ws.Cells.Item(137, 11).Value = 2579.3529  # K137
$ws.Cells.Item(137, 13).Value = -29.35289999999986  # M137
$ws.Cells.Item(138, 8).Value = 622457.8  # H138
$ws.Cells.Item(138, 9).Value = 905.75  # I138
$ws.Cells.Item(138, 10).Value = 1280571.9  # J138
$ws.Cells.Item(138, 11).Value = 2717.25  # K138
$ws.Cells.Item(138, 12).Value = 3841715.7  # L138
$ws.Cells.Item(138, 13).Value = 2422.75  # M138
$ws.Cells.Item(138, 14).Value = -3851995.7  # N138
$ws.Cells.Item(141, 8).Value = 494.3158  # H141
$ws.Cells.Item(141, 9).Value = 494.3158  # I141
$ws.Cells.Item(141, 11).Value = 1482.9474  # K141
$ws.Cells.Item(141, 13).Value = 3697.0526  # M141
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3328.7605  # H32
$ws.Cells.Item(32, 9).Value = 3099.1594  # I32
$ws.Cells.Item(32, 11).Value = 3099.1594  # K32
$ws.Cells.Item(32, 13).Value = -2812.1594  # M32
$ws.Cells.Item(61, 8).Value = 29412854  # H61
$ws.Cells.Item(61, 9).Value = 37038064  # I61
$ws.Cells.Item(61, 10).Value = 1328.7142  # J61
$ws.Cells.Item(61, 11).Value = 37038064  # K61
$ws.Cells.Item(61, 12).Value = 1328.7142  # L61
$ws.Cells.Item(61, 13).Value = -37037852  # M61
$ws.Cells.Item(61, 14).Value = -1752.7142  # N61
$ws.Cells.Item(74, 8).Value = 953.2093  # H74
$ws.Cells.Item(74, 9).Value = 670.3684  # I74
$ws.Cells.Item(74, 10).Value = 3102.8  # J74
$ws.Cells.Item(74, 11).Value = 670.3684  # K74
$ws.Cells.Item(74, 12).Value = 3102.8  # L74
$ws.Cells.Item(74, 13).Value = 203.6316  # M74
$ws.Cells.Item(74, 14).Value = -4850.8  # N74
$ws.Cells.Item(77, 8).Value = 953.2093  # H77
$ws.Cells.Item(77, 9).Value = 670.3684  # I77
$ws.Cells.Item(77, 10).Value = 3102.8  # J77
$ws.Cells.Item(77, 11).Value = 3351.842  # K77
$ws.Cells.Item(77, 12).Value = 15514  # L77
$ws.Cells.Item(77, 13).Value = 1016.158  # M77
$ws.Cells.Item(77, 14).Value = -24250  # N77
$ws.Cells.Item(102, 8).Value = 18519232  # H102
$ws.Cells.Item(102, 9).Value = 18519232  # I102
$ws.Cells.Item(102, 11).Value = 18519232  # K102
$ws.Cells.Item(102, 13).Value = -18517610  # M102
$ws.Cells.Item(132, 8).Value = 1459.3611  # H132
$ws.Cells.Item(132, 9).Value = 1351.8704  # I132
$ws.Cells.Item(132, 11).Value = 4055.6112  # K132
$ws.Cells.Item(132, 13).Value = -1525.6112  # M132
$ws.Cells.Item(133, 8).Value = 33200  # H133
$ws.Cells.Item(133, 10).Value = 33200  # J133
$ws.Cells.Item(133, 12).Value = 33200  # L133
$ws.Cells.Item(133, 14).Value = -38260  # N133
$ws.Cells.Item(136, 8).Value = 29412854  # H136
$ws.Cells.Item(136, 9).Value = 37038064  # I136
$ws.Cells.Item(136, 10).Value = 1328.7142  # J136
$ws.Cells.Item(136, 11).Value = 111114192  # K136
$ws.Cells.Item(136, 12).Value = 3986.1426  # L136
$ws.Cells.Item(136, 13).Value = -111111642  # M136
$ws.Cells.Item(136, 14).Value = -9086.142599999999  # N136
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3626.25  # H20
$ws.Cells.Item(20, 9).Value = 3500  # I20
$ws.Cells.Item(20, 11).Value = 3500  # K20
$ws.Cells.Item(20, 13).Value = -3253  # M20
$ws.Cells.Item(94, 8).Value = 125001500  # H94
$ws.Cells.Item(94, 9).Value = 250000000  # I94
$ws.Cells.Item(94, 10).Value = 3010  # J94
$ws.Cells.Item(94, 11).Value = 250000000  # K94
$ws.Cells.Item(94, 12).Value = 3010  # L94
$ws.Cells.Item(94, 13).Value = -249999549  # M94
$ws.Cells.Item(94, 14).Value = -3912  # N94
$ws.Cells.Item(99, 8).Value = 40000956  # H99
$ws.Cells.Item(99, 9).Value = 52632280  # I99
$ws.Cells.Item(99, 10).Value = 1766.5  # J99
$ws.Cells.Item(99, 11).Value = 52632280  # K99
$ws.Cells.Item(99, 12).Value = 1766.5  # L99
$ws.Cells.Item(99, 13).Value = -52630782  # M99
$ws.Cells.Item(99, 14).Value = -4762.5  # N99
$ws.Cells.Item(134, 8).Value = 2980.0625  # H134
$ws.Cells.Item(134, 9).Value = 929.56366  # I134
$ws.Cells.Item(134, 10).Value = 15510.889  # J134
$ws.Cells.Item(134, 11).Value = 2788.69098  # K134
$ws.Cells.Item(134, 12).Value = 46532.667  # L134
$ws.Cells.Item(134, 13).Value = -253.6909800000003  # M134
$ws.Cells.Item(134, 14).Value = -51602.667  # N134
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 100001544  # H16
$ws.Cells.Item(16, 9).Value = 111112670  # I16
$ws.Cells.Item(16, 10).Value = 1413  # J16
$ws.Cells.Item(16, 11).Value = 111112670  # K16
$ws.Cells.Item(16, 12).Value = 1413  # L16
$ws.Cells.Item(16, 13).Value = -111112383  # M16
$ws.Cells.Item(16, 14).Value = -1987  # N16
$ws.Cells.Item(31, 8).Value = 1867.5333  # H31
$ws.Cells.Item(31, 9).Value = 2007.3334  # I31
$ws.Cells.Item(31, 10).Value = 1308.3334  # J31
$ws.Cells.Item(31, 11).Value = 2007.3334  # K31
$ws.Cells.Item(31, 12).Value = 1308.3334  # L31
$ws.Cells.Item(31, 13).Value = -1712.3334  # M31
$ws.Cells.Item(31, 14).Value = -1898.3334  # N31
$ws.Cells.Item(34, 8).Value = 1867.5333  # H34
$ws.Cells.Item(34, 9).Value = 2007.3334  # I34
$ws.Cells.Item(34, 10).Value = 1308.3334  # J34
$ws.Cells.Item(34, 11).Value = 2007.3334  # K34
$ws.Cells.Item(34, 12).Value = 1308.3334  # L34
$ws.Cells.Item(34, 13).Value = -1805.3334  # M34
$ws.Cells.Item(34, 14).Value = -1712.3334  # N34
$ws.Cells.Item(111, 8).Value = 43499.5  # H111
$ws.Cells.Item(111, 10).Value = 43499.5  # J111
$ws.Cells.Item(111, 12).Value = 43499.5  # L111
$ws.Cells.Item(111, 14).Value = -51679.5  # N111
$ws.Cells.Item(113, 8).Value = 100001544  # H113
$ws.Cells.Item(113, 9).Value = 111112670  # I113
$ws.Cells.Item(113, 10).Value = 1413  # J113
$ws.Cells.Item(113, 11).Value = 111112670  # K113
$ws.Cells.Item(113, 12).Value = 1413  # L113
$ws.Cells.Item(113, 13).Value = -111110500  # M113
$ws.Cells.Item(113, 14).Value = -5753  # N113
$ws.Cells.Item(132, 8).Value = 3142.5425  # H132
$ws.Cells.Item(132, 9).Value = 3010.8628  # I132
$ws.Cells.Item(132, 11).Value = 9032.588400000001  # K132
$ws.Cells.Item(132, 13).Value = -6502.588400000001  # M132
$ws.Cells.Item(134, 8).Value = 8621653  # H134
$ws.Cells.Item(134, 9).Value = 994.59186  # I134
$ws.Cells.Item(134, 10).Value = 55556344  # J134
$ws.Cells.Item(134, 11).Value = 2983.77558  # K134
$ws.Cells.Item(134, 12).Value = 166669032  # L134
$ws.Cells.Item(134, 13).Value = -448.77558  # M134
$ws.Cells.Item(134, 14).Value = -166674102  # N134
$ws.Cells.Item(135, 8).Value = 34460  # H135
$ws.Cells.Item(135, 10).Value = 34460  # J135
$ws.Cells.Item(135, 12).Value = 34460  # L135
$ws.Cells.Item(135, 14).Value = -44600  # N135
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1829.8636  # H5
$ws.Cells.Item(5, 9).Value = 1985.7059  # I5
$ws.Cells.Item(5, 11).Value = 5957.1177  # K5
$ws.Cells.Item(5, 13).Value = -5845.1177  # M5
$ws.Cells.Item(80, 8).Value = 3445.4546  # H80
$ws.Cells.Item(80, 9).Value = 1133.3334  # I80
$ws.Cells.Item(80, 11).Value = 3400.0002  # K80
$ws.Cells.Item(80, 13).Value = -2464.0002  # M80
$ws.Cells.Item(83, 8).Value = 3445.4546  # H83
$ws.Cells.Item(83, 9).Value = 1133.3334  # I83
$ws.Cells.Item(83, 11).Value = 10200.0006  # K83
$ws.Cells.Item(83, 13).Value = -5520.000599999999  # M83
$ws.Cells.Item(112, 8).Value = 62511160  # H112
$ws.Cells.Item(112, 9).Value = 3513.25  # I112
$ws.Cells.Item(112, 10).Value = 83347040  # J112
$ws.Cells.Item(112, 11).Value = 10539.75  # K112
$ws.Cells.Item(112, 12).Value = 250041120  # L112
$ws.Cells.Item(112, 13).Value = -9431.75  # M112
$ws.Cells.Item(112, 14).Value = -250043336  # N112
$ws.Cells.Item(116, 8).Value = 3665  # H116
$ws.Cells.Item(116, 9).Value = 3325  # I116
$ws.Cells.Item(116, 10).Value = 3750  # J116
$ws.Cells.Item(116, 11).Value = 9975  # K116
$ws.Cells.Item(116, 12).Value = 11250  # L116
$ws.Cells.Item(116, 14).Value = -18134  # N116
$ws.Cells.Item(116, 13).Value = -6533  # M116
$ws.Cells.Item(120, 8).Value = 7005.6924  # H120
$ws.Cells.Item(120, 9).Value = 2009  # I120
$ws.Cells.Item(120, 10).Value = 10128.625  # J120
$ws.Cells.Item(120, 11).Value = 6027  # K120
$ws.Cells.Item(120, 12).Value = 30385.875  # L120
$ws.Cells.Item(120, 13).Value = -1189  # M120
$ws.Cells.Item(120, 14).Value = -40061.875  # N120
$ws.Cells.Item(135, 8).Value = 1829.8636  # H135
$ws.Cells.Item(135, 9).Value = 1985.7059  # I135
$ws.Cells.Item(135, 11).Value = 17871.3531  # K135
$ws.Cells.Item(135, 13).Value = -15336.3531  # M135
$ws.Cells.Item(139, 8).Value = 1576.5  # H139
$ws.Cells.Item(139, 9).Value = 1576.5  # I139
$ws.Cells.Item(139, 10).Value = 0  # J139
$ws.Cells.Item(139, 11).Value = 4729.5  # K139
$ws.Cells.Item(139, 12).Value = 0  # L139
$ws.Cells.Item(139, 13).Value = 410.5  # M139
$ws.Cells.Item(139, 14).ClearContents()  # N139
$ws.Cells.Item(140, 8).Value = 22765.64  # H140
$ws.Cells.Item(140, 9).Value = 64818.062  # I140
$ws.Cells.Item(140, 10).Value = 2976.2646  # J140
$ws.Cells.Item(140, 11).Value = 194454.186  # K140
$ws.Cells.Item(140, 12).Value = 8928.793799999999  # L140
$ws.Cells.Item(140, 13).Value = -189274.186  # M140
$ws.Cells.Item(140, 14).Value = -19288.7938  # N140
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 150001400  # H70
$ws.Cells.Item(70, 9).Value = 250000000  # I70
$ws.Cells.Item(70, 10).Value = 100002100  # J70
$ws.Cells.Item(70, 11).Value = 250000000  # K70
$ws.Cells.Item(70, 12).Value = 100002100  # L70
$ws.Cells.Item(70, 13).Value = -249999730  # M70
$ws.Cells.Item(70, 14).Value = -100002640  # N70
$ws.Cells.Item(73, 8).Value = 150001400  # H73
$ws.Cells.Item(73, 9).Value = 250000000  # I73
$ws.Cells.Item(73, 10).Value = 100002100  # J73
$ws.Cells.Item(73, 11).Value = 250000000  # K73
$ws.Cells.Item(73, 12).Value = 100002100  # L73
$ws.Cells.Item(73, 13).Value = -249999064  # M73
$ws.Cells.Item(73, 14).Value = -100003972  # N73
$ws.Cells.Item(97, 8).Value = 479.15384  # H97
$ws.Cells.Item(97, 9).Value = 479.15384  # I97
$ws.Cells.Item(97, 11).Value = 479.15384  # K97
$ws.Cells.Item(97, 13).Value = 16.84616  # M97
$ws.Cells.Item(132, 8).Value = 1366.9464  # H132
$ws.Cells.Item(132, 9).Value = 1196.4286  # I132
$ws.Cells.Item(132, 11).Value = 3589.2858  # K132
$ws.Cells.Item(132, 13).Value = -1059.2858  # M132
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(92, 8).Value = 15500  # H92
$ws.Cells.Item(92, 10).Value = 15500  # J92
$ws.Cells.Item(92, 12).Value = 15500  # L92
$ws.Cells.Item(92, 14).Value = -20492  # N92
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 4500  # H96
$ws.Cells.Item(96, 9).Value = 4500  # I96
$ws.Cells.Item(96, 10).Value = 0  # J96
$ws.Cells.Item(96, 11).Value = 4500  # K96
$ws.Cells.Item(96, 12).Value = 0  # L96
$ws.Cells.Item(96, 13).Value = -3127  # M96
$ws.Cells.Item(96, 14).ClearContents()  # N96
$ws.Cells.Item(100, 8).Value = 704.6667  # H100
$ws.Cells.Item(100, 10).Value = 646  # J100
$ws.Cells.Item(100, 12).Value = 1292  # L100
$ws.Cells.Item(100, 14).Value = -2374  # N100
$ws.Cells.Item(132, 8).Value = 2579.1025  # H132
$ws.Cells.Item(132, 9).Value = 2642.0605  # I132
$ws.Cells.Item(132, 11).Value = 7926.181500000001  # K132
$ws.Cells.Item(132, 13).Value = -5396.181500000001  # M132
$ws.Cells.Item(136, 8).Value = 539.8570999999999  # H136
$ws.Cells.Item(136, 9).Value = 375.5  # I136
$ws.Cells.Item(136, 10).Value = 868.5714  # J136
$ws.Cells.Item(136, 11).Value = 1126.5  # K136
$ws.Cells.Item(136, 12).Value = 2605.7142  # L136
$ws.Cells.Item(136, 13).Value = 1423.5  # M136
$ws.Cells.Item(136, 14).Value = -7705.7142  # N136
$ws.Cells.Item(141, 8).Value = 25070.889  # H141
$ws.Cells.Item(141, 10).Value = 33281.332  # J141
$ws.Cells.Item(141, 12).Value = 33281.332  # L141
$ws.Cells.Item(141, 14).Value = -43641.332  # N141
